# Report 18 02 2025
# Update the previously-recorded timestamp (rows 284-323, column A) with the
# corrected precision value, then append the newly scraped silver price rows
# (324-345) collected on 2025-02-18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = "YYYY-MM-DD HH:MM:SS"

# --- 1. Fix the slight floating point drift on the existing timestamp column ---
$fixedTimestamp = 45704.9909262963
for ($r = 284; $r -le 323; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $fixedTimestamp
    $ws.Cells.Item($r, 1).NumberFormat = $dateFormat
}

# --- 2. Append the new rows scraped on 2025-02-18 ---
$newTimestamp = 45706.87791302231

$newRows = @(
    @("1 килограм сребърно кюлче Valcambi", 2067.47, 2639.8, "https://tavex.bg/srebro/1-kg-valcambi-sreburno-kiulche/"),
    @("1/10 унция сребърна монета Британия", 9.23, 12.33, "https://tavex.bg/srebro/1-10-unciya-srebarna-moneta-britaniya/"),
    @("1 унция сребърна монета Виенска Филхармония", 64.61, 81.37, "https://tavex.bg/srebro/1-unciq-srebyrna-avstriiska-filharmonia/"),
    @("1 унция сребърна монета канадски кленов лист", 65.84, 84.31999999999999, "https://tavex.bg/srebro/1-unciya-sreburen-kanadski-klenov-list/"),
    @("1 унция сребърна монета Австралийско кенгуру", 64.61, 82.09999999999999, "https://tavex.bg/srebro/1-oz-sreburna-moneta-avstraliysko-kenguru/"),
    @("1 унция сребърна монета Американски орел", 67.69, 118.36, "https://tavex.bg/srebro/1-oz-sreburen-orel/"),
    @("1 унция сребърен австралийски лунар Змия 2025", 80, 177.53, "https://tavex.bg/srebro/1-unciya-srebaren-avstraliyski-lunar-godina-na-zmiyata-2025/"),
    @("30 грама сребърна монета Китайска панда 2025", 65.29000000000001, 156.96, "https://tavex.bg/srebro/30-grama-srebarna-moneta-kitaiska-panda-2025/"),
    @("1 унция сребърен австралийски лунар Дракон 2024", 73.84, 155.34, "https://tavex.bg/srebro/1-unciya-srebyren-avstraliiski-lunar-drakon-2024/"),
    @("1 унция сребърен австралийски лунар Заек 2023", 80, 207.12, "https://tavex.bg/srebro/1-unciya-srebyren-avstraliiski-lunar-zaek-2023/"),
    @("30 грама сребърна монета Китайска панда 2024", 80.29000000000001, $null, "https://tavex.bg/srebro/30-grama-srebarna-moneta-kitayska-panda-2024/"),
    @("30 грама сребърна монета Китайска панда 2023", 80.29000000000001, $null, "https://tavex.bg/srebro/30-grama-srebarna-moneta-kitayska-panda-2023/"),
    @("1 унция сребърна монета Британия", 80.84, $null, "https://tavex.bg/srebro/1-unciya-srebarna-moneta-britaniya-2/"),
    @("1 унция Сребърна монета Кругерранд, Южна Африка", 80.84, $null, "https://tavex.bg/srebro/1-unciya-srebarna-moneta-krugerrand-yuzhna-afrika/"),
    @("25 бр. 1 унция сребърна монета Британия", 80, $null, "https://tavex.bg/srebro/25-broya-1-unciya-srebarna-moneta-britania-tubus/"),
    @("25 бр. 1 унция сребърна монета Кругерранд, Южна Африка", 80, $null, "https://tavex.bg/srebro/25-broya-1-unciya-srebarna-moneta-krugerrand-yujna-afrika/"),
    @("500 бр. 1 унция Сребърна монета Британия", 80, $null, "https://tavex.bg/srebro/500-broya-1-unciya-srebarna-moneta-britaniya-masterbox-kutiya/"),
    @("500 бр. 1 унция сребърна монета Кругерранд, Южна Африка", 80, $null, "https://tavex.bg/srebro/500-broya-1-unciya-srebarna-moneta-krugerrand-yujna-afrika/"),
    @("1 унция сребърна австралийска коала", 80.69, $null, "https://tavex.bg/srebro/1-unciya-srebyrna-avstraliiska-koala/"),
    @("1 унция сребърна монета австралийски лунар година на Тигъра 2022", 80.69, $null, "https://tavex.bg/srebro/1-unciya-srebyrna-moneta-avstraliiski-lunar-tigyr-2023/"),
    @("30 грама сребърна монета Китайска панда 2022", 80.29000000000001, $null, "https://tavex.bg/srebro/30-grama-srebyrna-kitayska-panda-2022/"),
    @("1 унция  Кукабура 2022 година", 80.69, $null, "https://tavex.bg/srebro/1-oz-australian-kookaburra-2022-silver-coin/")
)

$startRow = 324
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value2 = $newTimestamp
    $ws.Cells.Item($r, 1).NumberFormat = $dateFormat

    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]

    if ($null -eq $row[2]) {
        # Keep the cell present but empty (matches an inline blank cell in the
        # source data) instead of leaving it completely absent from the sheet.
        $ws.Cells.Item($r, 4).Value = ""
        $ws.Cells.Item($r, 4).Style = "Normal"
    } else {
        $ws.Cells.Item($r, 4).Value = $row[2]
    }

    $ws.Cells.Item($r, 5).Value = $row[3]
}
